$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2021 active list")

# ---- Header row (row 1) renames ----
$ws.Range("A1").Value = "#"
$ws.Range("G1").Value = "2019 Total"
$ws.Range("H1").Value = "Returning #"
$ws.Range("I1").Value = "2020 Total"
$ws.Range("J1").Value = "Max # 2021"
$ws.Range("K1").Value = "Requested/ confirmed #"
$ws.Range("L1").Value = "Coordinator"
$ws.Range("O1").Value = "Passport presentation"
$ws.Range("P1").Value = "Portal / Passports"
$ws.Range("Q1").Value = "Agreement/ consent"
$ws.Range("R1").Value = "Notes"

# ---- Row 2 (Amuri Area School) data updates ----
$ws.Range("G2").Value = 19
$ws.Range("H2").Value = 8
$ws.Range("I2").Value = 35
$ws.Range("P2").Value = "Y/Y"
$ws.Range("Q2").Value = "Y/Y"
$ws.Range("R2").Value = "Agreement signed, consent form recieved"

# ---- Row 3 (Rawhiti School) data updates ----
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = "NA"
$ws.Range("I3").Value = 25
$ws.Range("P3").Value = "Y/Y"
$ws.Range("Q3").Value = "Y/Y"
$ws.Range("R3").Value = "Agreement signed"

# ---- Re-point the cell formatting to match the re-shuffled columns:
#      "2019 Total" / row totals now use the centered numeric look, while
#      the merged "Notes" column reuses the plain left/vert-centered look ----
$ws.Range("H1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

$ws.Range("A1").Copy()
$ws.Range("R1").PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Range("G2").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("R2").PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Range("G3").PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Range("P3").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$ws.Range("R3").PasteSpecial(-4122)

# ---- Remove the stray empty styled cell that used to sit under the old
#      "Year" column, keeping the note text in F9 intact (with its style) ----
$ws.Rows(9).Delete()
$ws.Range("A3").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F9").Value = "In Progress = waiting on paperwork for 2021"

# ---- The old "Notes" column (T) has been folded into the "Agreement" column
#      (now renamed "Notes", column R above); the old stand-alone "Consent"
#      column (S) is no longer needed. Remove both trailing columns. ----
$ws.Range("T1:T9").EntireColumn.Delete()
$ws.Range("S1:S9").EntireColumn.Delete()
